$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the data range to Text format before writing, so that numeric-looking
# strings (e.g. "624.10", "0.607") are preserved exactly as text instead of being
# auto-converted to floating point numbers by Excel. Revert the style afterwards
# so cells keep their original (default) style.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "70.602.95"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "3.510.28"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "624.10"
$ws.Range("E5").Value = "  +4.10%  "
$ws.Range("D6").Value = "171.73"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").Value = "3.505.35"
$ws.Range("E7").Value = "  -1.77%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "0.607"
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").Value = "7.18"
$ws.Range("E11").Value = "  -2.62%  "
$ws.Range("D12").Value = "0.583"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").Value = "46.13"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").Value = "0.0000275"
$ws.Range("E14").Value = "  -1.12%  "
$ws.Range("D15").Value = "4.078.65"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("D16").Value = "8.41"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "604.54"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").Value = "3.518.26"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").Value = "70.745.78"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("D21").Value = "17.64"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").Value = "0.877"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").Value = "9.07"
$ws.Range("E23").Value = "  -1.28%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "96.96"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "15.45"
$ws.Range("E25").Value = "  -2.66%  "
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  -3.53%  "
$ws.Range("D29").Value = "33.39"
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").Value = "8.98"
$ws.Range("E30").Value = "  -2.10%  "
$ws.Range("D31").Value = "2.99"
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("D32").Value = "8.05"
$ws.Range("E32").Value = "  -3.73%  "
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("D34").Value = "6.77"
$ws.Range("E34").Value = "  -5.84%  "
$ws.Range("D35").Value = "621.63"
$ws.Range("E35").Value = "  -6.88%  "
$ws.Range("B36").Value = "Cosmos"
$ws.Range("C36").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D36").Value = "10.81"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.0489"
$ws.Range("E37").Value = "  +2.04%  "
$ws.Range("D38").Value = "0.0991"
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("D39").Value = "3.40"
$ws.Range("E39").Value = "  -7.86%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "56.62"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").Value = "3.328.64"
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("D44").Value = "0.0₃0720"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").Value = "0.309"
$ws.Range("E46").Value = "  -4.29%  "
$ws.Range("D47").Value = "31.73"
$ws.Range("E47").Value = "  -3.74%  "
$ws.Range("D48").Value = "2.48"
$ws.Range("E48").Value = "  -6.90%  "
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").Value = "133.92"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("E51").Value = "  +0.00%  "

# Restore the original (default) cell style now that values are safely stored as text.
$dataRange.Style = "Normal"
